$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.893.12"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "1.840.36"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.34"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.617"
$ws.Range("E6").Value = "  +2.47%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.91"
$ws.Range("E8").Value = "  +4.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.327"
$ws.Range("E9").Value = "  +3.49%  "
$ws.Range("E10").Value = "  +1.88%  "
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").Value = "2.106.87"
$ws.Range("E12").Value = "  +1.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.33"
$ws.Range("E13").Value = "  +4.47%  "
$ws.Range("D14").Value = "1.840.63"
$ws.Range("E14").Value = "  +1.90%  "
$ws.Range("E15").Value = "  +1.67%  "
$ws.Range("E16").Value = "  +2.58%  "
$ws.Range("D17").Value = "34.942.51"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.80"
$ws.Range("E19").Value = "  +1.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "239.84"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("E21").Value = "  +3.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.74"
$ws.Range("E22").Value = "  +2.34%  "
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").Value = "  +1.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.02"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("E26").Value = "  +2.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.38"
$ws.Range("E27").Value = "  +1.67%  "
$ws.Range("E28").Value = "  +3.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.65"
$ws.Range("E29").Value = "  +9.00%  "
$ws.Range("E31").Value = "  +1.51%  "
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.90"
$ws.Range("E33").Value = "  -0.64%  "
$ws.Range("E34").Value = "  +22.52%  "
$ws.Range("E35").Value = "  +10.62%  "
$ws.Range("E36").Value = "  +0.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.745"
$ws.Range("E37").Value = "  +8.23%  "
$ws.Range("E38").Value = "  +10.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "89.68"
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("E40").Value = "  +3.55%  "
$ws.Range("D41").Value = "1.337.12"
$ws.Range("E41").Value = "  +2.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.55"
$ws.Range("E42").Value = "  +2.89%  "
$ws.Range("E43").Value = "  -1.94%  "
$ws.Range("E44").Value = "  +2.06%  "
$ws.Range("E45").Value = "  +3.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0531"
$ws.Range("E46").Value = "  +4.15%  "
$ws.Range("E47").Value = "  +3.09%  "
$ws.Range("D48").Value = "2.025.66"
$ws.Range("E48").Value = "  +1.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.72"
$ws.Range("E49").Value = "  +66.99%  "
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0667"
$ws.Range("E51").Value = "  -0.42%  "
